# Generate Report for Handback
# The handback transform failed for file 61bd1201-727f-4028-9b78-6173754d0f3c
# (row 7 of the zh-cn / de-de sheets, and the corresponding row of Overview).
# Update status text everywhere it is shown, and record the error detail in
# column L ("Error Detail") of the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the "Ready for handoff" status to "Handback transform failed"
# everywhere it appears (Overview zh-cn/de-de columns, and the Status
# column on each language detail sheet) for the affected file's row.
$wsOverview.Range("B7").Value = $newStatus
$wsOverview.Range("C7").Value = $newStatus
$wsZhCn.Range("C7").Value = $newStatus
$wsDeDe.Range("C7").Value = $newStatus

# Record the error detail explaining the handback transform failure.
$wsZhCn.Range("L7").Value = "Handback file name: koowryoq.21w is different with handoff file name: 61bd1201-727f-4028-9b78-6173754d0f3c.deda338c481ed7b85ea955ef1b2db1f48ebf6515.zh-cn."
$wsDeDe.Range("L7").Value = "Handback file name: koowryoq.21w is different with handoff file name: 61bd1201-727f-4028-9b78-6173754d0f3c.deda338c481ed7b85ea955ef1b2db1f48ebf6515.de-de."
